$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Real Madrid CF - Cadiz CF", "19/12/2021"),
    @("Real Madrid CF - Valencia CF", "09/01/2022"),
    @("Real Madrid CF - Elche CF", "23/01/2022"),
    @("Real Madrid CF - Granada CF", "06/02/2022"),
    @("Real Madrid CF - Deportivo Alavés", "20/02/2022"),
    @("Real Madrid CF - Real Sociedad", "06/03/2022"),
    @("Real Madrid CF - FC Barcelona", "20/03/2022"),
    @("Real Madrid CF - Getafe CF", "10/04/2022"),
    @("Real Madrid CF - RCD Espanyol", "01/05/2022"),
    @("Real Madrid CF - UD Levante", "11/05/2022")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    $cellA.Value2 = $data[$i][0]

    # Column B holds DD/MM/YYYY-style text. Some of these (day <= 12) are
    # ambiguous and would otherwise get silently reinterpreted as a date
    # serial number. Force text entry, then restore the default "Normal"
    # style so the cell is left exactly as it was (plain shared-string text,
    # default formatting).
    $cellB.NumberFormat = "@"
    $cellB.Value2 = $data[$i][1]
    $cellB.Style = "Normal"
}
